$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror column J formatting into the new column K (year 2020 / value 173)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2020

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 173

# Update selection on the sheet view
$ws.Range("I18").Select()
